$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-18 Thursday" "2025-12-19 Friday"

Replace-Text "81×37=2997" "42×82=3444"
Replace-Text "85×45=3825" "35×73=2555"
Replace-Text "26×61=1586" "90×52=4680"
Replace-Text "35×91=3185" "36×74=2664"
Replace-Text "39×88=3432" "23×62=1426"

Replace-Text "51×55=2805" "81×25=2025"
Replace-Text "76×94=7144" "95×67=6365"
Replace-Text "39×16=624" "44×84=3696"
Replace-Text "35×36=1260" "74×79=5846"
Replace-Text "91×34=3094" "56×55=3080"

Replace-Text "48×26=1248" "59×91=5369"
Replace-Text "18×70=1260" "87×91=7917"
Replace-Text "52×78=4056" "94×18=1692"
Replace-Text "92×30=2760" "26×23=598"
Replace-Text "51×39=1989" "13×12=156"

Replace-Text "43×67=2881" "96×26=2496"
Replace-Text "95×74=7030" "30×93=2790"
Replace-Text "64×89=5696" "72×97=6984"
Replace-Text "26×90=2340" "12×49=588"
Replace-Text "96×99=9504" "55×37=2035"

Replace-Text "96×23=2208" "38×93=3534"
Replace-Text "21×79=1659" "14×94=1316"
Replace-Text "72×73=5256" "29×73=2117"
Replace-Text "69×94=6486" "58×70=4060"
Replace-Text "89×52=4628" "89×37=3293"
